$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 212, shifting existing rows 212:304 down to 213:305.
# The new row inherits formatting (incl. the date number format on column D) from
# the row above it, matching the original workbook's layout.
$ws.Rows("212:212").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A212").Value = 8
$ws.Range("B212").Value = "Terminal La Palmera de La Serena"
$ws.Range("C212").Value = "Coquimbo"
$ws.Range("D212").Value = 45202
$ws.Range("E212").Value = 4
$ws.Range("F212").Value = 100112001
$ws.Range("G212").Value = "Berenjena"
$ws.Range("H212").Value = "Sin especificar"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 560
$ws.Range("K212").Value = 9000
$ws.Range("L212").Value = 10000
$ws.Range("M212").Value = 9500
$ws.Range("N212").Value = "$/caja 50 unidades"
$ws.Range("O212").Value = "Región de Arica y Parinacota"
$ws.Range("P212").Value = 190
$ws.Range("Q212").Value = 50
$ws.Range("R212").Value = "Hortaliza"
